# Scottish Module Input files modification
#
# This payroll input workbook was repurposed from the "Monthly" tax
# variant to the "Weekly" tax variant:
#   - the "GeneralTaxRateMonthly" sheet becomes "GeneralTaxRateWeekly"
#   - the "ProcessPayrollForMonthlyTax" sheet becomes "ProcessPayrollForWeeklyTax"
#   - the overview sheet ("first") is updated to reference the new names
#   - the "DO NOT TOUCH AUTOMATION EMP 105" marker becomes
#     "DO NOT TOUCH AUTOMATION EMP 107" on every sheet that carries it

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("first")
$ws2 = $wb.Worksheets.Item("GeneralTaxRateMonthly")
$ws3 = $wb.Worksheets.Item("ProcessPayrollForMonthlyTax")
$ws4 = $wb.Worksheets.Item("TestReports")

# Bump the employee automation marker on each sheet that references it
# (done before the renames below so the shared-string table keeps its
# natural insertion order).
$ws2.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$ws3.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$ws4.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# Rename the monthly-tax sheets to their weekly-tax equivalents
$ws2.Name = "GeneralTaxRateWeekly"
$ws3.Name = "ProcessPayrollForWeeklyTax"

# Keep the "first" overview sheet's references in sync with the renamed sheets
$ws1.Range("A3").Value = "GeneralTaxRateWeekly"
$ws1.Range("A4").Value = "ProcessPayrollForWeeklyTax"

# Leave the workbook with the same view/selection state it was saved with:
# "first" ends up the active sheet (selection A3), while the other sheets
# keep whatever cell was last selected on them.
$ws2.Activate()
$ws2.Range("E15").Select()

$ws3.Activate()
$ws3.Range("B2").Select()

$ws4.Activate()
$ws4.Range("B2").Select()

$ws1.Activate()
$ws1.Range("A3").Select()
